$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) uses a custom date/time number format throughout the
# sheet (the same format as every other data row) - grab it once so the
# rows we touch/add keep matching formatting instead of falling back to
# the generic "General" number format.
$dateFmt = $ws.Cells.Item(70, 4).NumberFormat

# ---------------------------------------------------------------------------
# The sheet held 72 data rows (rows 2-73). This edit:
#   1. Updates existing row 71 with new values (new date/variety/prices/etc).
#   2. Updates existing row 72 with new values (same date as row 71 now).
#   3. Replaces row 73 with a brand-new "Inferno" record (new date).
#   4. Appends three new rows (74-76) that preserve the ORIGINAL data that
#      used to live in rows 71, 72 and 73 before this edit.
# ---------------------------------------------------------------------------

function Set-DataRow($Row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 4).NumberFormat = $dateFmt
}

# Row 71 - updated: Americana(o), 25kg box, Provincia de Limarí
Set-DataRow 71 11 "Vega Monumental Concepción" "Bíobío" 44610 8 100112021 `
    "Ají" "Americana (o)" "Primera" 100 17000 18000 17500 `
    "$/caja 25 kilos" "Provincia de Limarí" 700 25 "Hortaliza"

# Row 72 - updated: Chilena(o), 25kg box, Región Metropolitana
Set-DataRow 72 11 "Vega Monumental Concepción" "Bíobío" 44610 8 100112021 `
    "Ají" "Chilena(o)" "Primera" 50 25000 26000 25400 `
    "$/caja 25 kilos" "Región Metropolitana" 1016 25 "Hortaliza"

# Row 73 - new record: Inferno, 12kg box, Región de Arica y Parinacota
Set-DataRow 73 11 "Vega Monumental Concepción" "Bíobío" 44610 8 100112021 `
    "Ají" "Inferno" "Primera" 50 12000 13000 12600 `
    "$/caja 12 kilos" "Región de Arica y Parinacota" 1050 12 "Hortaliza"

# Row 74 - appended: original data previously held in row 71
Set-DataRow 74 11 "Vega Monumental Concepción" "Bíobío" 44399 8 100112021 `
    "Ají" "Inferno" "Primera" 50 25000 26000 25600 `
    "$/caja 12 kilos" "Región de Arica y Parinacota" 2133 12 "Hortaliza"

# Row 75 - appended: original data previously held in row 72
Set-DataRow 75 11 "Vega Monumental Concepción" "Bíobío" 44595 8 100112021 `
    "Ají" "Chilena(o)" "Primera" 100 22000 23000 22500 `
    "$/caja 25 kilos" "Región Metropolitana" 900 25 "Hortaliza"

# Row 76 - appended: original data previously held in row 73
Set-DataRow 76 11 "Vega Monumental Concepción" "Bíobío" 44552 8 100112021 `
    "Ají" "Americana (o)" "Primera" 40 36000 38000 37000 `
    "$/caja 25 kilos" "Provincia de Limarí" 1480 25 "Hortaliza"
